$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5082618594169617
$ws.Range("B1").Value = 1.803717970848083
$ws.Range("C1").Value = 5.891898632049561
$ws.Range("D1").Value = 1.581515431404114
$ws.Range("E1").Value = 0.8052790760993958
